$wb = $excel.ActiveWorkbook

# --- Sheet "Canada": add new data row for date 2022-01-01 (serial 44562) ---
$wsCanada = $wb.Worksheets.Item("Canada")
$wsCanada.Range("A26").Value = 44562
$wsCanada.Range("B26").Value = "Canada"
$wsCanada.Range("C26").Formula = "=(D26-E26)/E26*100"
$wsCanada.Range("D26").Value = 1341.8
$wsCanada.Range("E26").Value = 1180.9000000000001
$wsCanada.Range("A26").NumberFormat = $wsCanada.Range("A25").NumberFormat
$wsCanada.Range("B26").NumberFormat = $wsCanada.Range("B25").NumberFormat

# --- Sheet "Province": add 10 new data rows for date 2022-01-01 (serial 44562) ---
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
  @{Row=242; Name="Newfoundland & Labrador"; D=32.2;               E=31.3},
  @{Row=243; Name="Prince Edward Island";    D=8.6;                E=8.1999999999999993},
  @{Row=244; Name="Nova Scotia";             D=35.799999999999997; E=35.1},
  @{Row=245; Name="New Brunswick";           D=33.4;               E=32.299999999999997},
  @{Row=246; Name="Quebec";                  D=247.3;              E=242.2},
  @{Row=247; Name="Ontario";                 D=591.4;              E=456.4},
  @{Row=248; Name="Manitoba";                D=35.5;               E=40.200000000000003},
  @{Row=249; Name="Saskatchewan";            D=33.299999999999997; E=34.1},
  @{Row=250; Name="Alberta";                 D=179.9;              E=167},
  @{Row=251; Name="British Columbia";        D=144.30000000000001; E=134.1}
)

foreach ($r in $provinceRows) {
  $row = $r.Row
  $wsProvince.Range("A$row").Value = 44562
  $wsProvince.Range("B$row").Value = $r.Name
  $wsProvince.Range("C$row").Formula = "=(D$row-E$row)/E$row*100"
  $wsProvince.Range("D$row").Value = $r.D
  $wsProvince.Range("E$row").Value = $r.E
  $wsProvince.Range("A$row").NumberFormat = $wsProvince.Range("A241").NumberFormat
}
$wsProvince.Range("B242").NumberFormat = $wsProvince.Range("B232").NumberFormat

# --- Update view/selection state to match the new data extents ---
# Visit "Canada" first and leave its selection on the newly added row.
[void]$wsCanada.Activate()
[void]$wsCanada.Range("A26").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1

# Finish on "Province" (the sheet that remains active in the workbook).
[void]$wsProvince.Activate()
[void]$wsProvince.Range("D252").Select()
$excel.ActiveWindow.ScrollRow = 238
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
